# Apply the "cryptos list" refresh described by the commit diff.
# Columns B/C (Coin name / Link) and D (Price) are plain/numeric-looking
# text that must remain TEXT (never auto-converted to a number by Excel),
# while column E (Volume/1h) already keeps its padding spaces so it stays
# text naturally.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.796.41'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +5.87%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.119.47'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.77%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '559.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.58%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.20'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +10.93%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.06%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.115.16'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.73%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.502'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.30%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.88'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +16.89%  '

# Row 11
$ws.Range("E11").Value = '  +6.49%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.462'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.44%  '

# Row 13
$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '35.99'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.23%  '

# Row 14
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000230'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.34%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.559.08'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.12%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.766.58'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.90%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.105.05'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.59%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.109'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.58%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.83'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.40%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '487.76'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.34%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.94'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.46%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.681'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.95%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.58%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.24'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +11.45%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.69'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.86%  '

# Row 26
$ws.Range("E26").Value = '  +0.07%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.82'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.17%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.08'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.86%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.09'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +11.25%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.15%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.52'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.78%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.16'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.20%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.48'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.89%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.83'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.26%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '55.62'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.10%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.16'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.08%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '468.06'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.20%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0838'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.02%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0408'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +9.01%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.030.05'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.20%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.119'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.39%  '

# Row 42
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.76'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +17.89%  '

# Row 43
$ws.Range("B43").Value = 'Cosmos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.36'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.30%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '28.18'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +10.93%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.262'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +9.88%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.09'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +8.83%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.112'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.00%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₃0523'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.38%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '117.60'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.28%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.10'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.42%  '
